# Update column F ("dSF") values on the active sheet for the rows whose
# value changed (repulled data / recalculated means).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 5
    6  = 2
    8  = -3
    9  = 9
    10 = -2
    11 = 1
    12 = 1
    13 = -8
    14 = 4
    15 = -5
    16 = -1
    17 = 6
    18 = 7
    19 = -5
    20 = 4
    21 = -4
    22 = -1
    23 = 8
    24 = -1
    25 = 1
    26 = -1
    27 = 6
    28 = -3
    29 = 4
    30 = 1
    31 = 2
    32 = -1
    33 = 6
    34 = -3
    36 = -1
    37 = -2
    38 = 5
    39 = 2
    40 = -1
    41 = -2
    42 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
